$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-21 10:43:44"
$wsZh.Range("H2").Value = "2016-03-21 10:44:05"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-21 10:43:47"
$wsDe.Range("H2").Value = "2016-03-21 10:44:11"
